$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44326, 4, 17, 168.8182720953327),
    @(44327, 1, 18, 178.7487586891758),
    @(44328, 0, 18, 178.7487586891758),
    @(44329, 0, 13, 129.0963257199603)
)

$startRow = 252
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $values = $data[$i]

    # Copy the format of the last existing date cell (A251, style "2": centered,
    # bordered, YYYY-MM-DD HH:MM:SS) onto the new date cell so the same style is reused.
    $ws.Range("A251").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value2 = $values[0]
    $ws.Cells.Item($r, 2).Value2 = $values[1]
    $ws.Cells.Item($r, 3).Value2 = $values[2]
    $ws.Cells.Item($r, 4).Value2 = $values[3]
}
